$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "xyz"
$ws.Range("A20").Value = "xyz"
